$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Row 3
$ws.Range("D3").Value = 10.19
$ws.Range("E3").Value = 10.25

# Row 4
$ws.Range("C4").Value = 9.81
$ws.Range("E4").Value = 9.890000000000001

# Row 5
$ws.Range("C5").Value = 9.75
$ws.Range("D5").Value = 10.11
$ws.Range("F5").Value = 10.1
$ws.Range("G5").Value = 9.4
$ws.Range("H5").Value = 8.130000000000001

# Row 6
$ws.Range("E6").Value = 9.9
$ws.Range("G6").Value = 10.15
$ws.Range("H6").Value = 11.91

# Row 7
$ws.Range("E7").Value = 10.6
$ws.Range("F7").Value = 9.85

# Row 8
$ws.Range("E8").Value = 11.87
$ws.Range("F8").Value = 8.09
